$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap D2/E2 values ("Different"/"Similar") before the column shift so the
# swap isn't affected by the later column deletion (columns D/E are unaffected
# by deleting column R, but do this first for clarity).
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$ws.Range("D2").Value2 = $e2
$ws.Range("E2").Value2 = $d2

# Remove the now-unused "block_length_multiplier" column (R) entirely,
# shifting columns S:Y left by one.
$ws.Range("R1").EntireColumn.Delete()

# Update the view: scroll the window so column B is the left-most visible
# column, then set the active selection to E3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E3").Select() | Out-Null
